$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    24 = "0148032060184c38150a70dc4cc62862"
    79 = "1ad46692ff209a26f17ea38a80419558"
    136 = "145f6cdd9e574970a49058607a4c57c6"
    150 = "19052ad2a734e1844672a4e4de748779"
    159 = "dbfc21f7e94c2499a7e91e097f364003"
    169 = "d8e2d3b430620fbcc36650018a5d213d"
    227 = "79d7ac27c02b8ee4b146a8ebaf9cdac1"
    232 = "ae22bcdb5a3d16e8e1bb7667b80435a8"
    281 = "7f6ab24a2600337270ff3e0396ae3efd"
    302 = "0f1ef506e706195dbd93c49065f789b1"
    339 = "4355b8ccd9f3d91560badc347230afcd"
    419 = "afba4ee92bb44bede48ddf483ac24705"
    460 = "ef3bb11c9a11290215fab20c3653025e"
    478 = "19b25a4ce25f6f97839a85d363ab88b0"
    500 = "90638a5840cb2ea45547ac598d99705e"
    501 = "10add39a694426657601535a2ecb2c04"
    502 = "81629ac93065ab0b8af54b4a0aeeec72"
    517 = "d58681c86cbed19c395aab18d70338ab"
    550 = "8aab137630c87b0adee966d8555f7e13"
    563 = "e36dde274970a017fcdcb0f19f6bba4c"
    616 = "078638d89707ef761041c1aa1f6eb798"
    627 = "0225aa8685f6b6a513936ce0d53587e9"
    665 = "1ba24c61578dfbe6dd75691af4a3de32"
    680 = "dfc9b3ba408aa959d34138ce25d08e59"
    685 = "5bc27490b7119c501eed5d24ed5b0700"
    700 = "c1be0d083ce0ad19eb1f14e63dd5771f"
    703 = "19cfb9675ed60fea946f78fdbeb00be0"
    704 = "aabab779119ff412ea1a22091217f45a"
    742 = "751a358c7da9cc56db9b1cdcade1cc19"
    819 = "ddcecae74f700d34aeb688e4eafe9966"
    830 = "878f501c6fcfbb24100b756563e49341"
    835 = "44a1dc031076aedec8ddf2465a2c79d5"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $updates[$row]
}

$wb.Save()
